# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F value
$updates = @{
    2  = 298
    4  = 10274
    5  = 329
    6  = 936
    7  = 1275
    8  = 6648
    10 = 434
    12 = 126
    13 = 3160
    15 = 309
    16 = 637
    18 = 603
    20 = 56
    21 = 1607
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
